$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20
$data[0,0] = "FAPs"
$data[0,1] = "Ntf5"
$data[0,2] = "Ngfr"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.113088
$data[0,7] = 0.339264
$data[0,8] = 0.3269336956678857
$data[0,9] = 0.3269336956678857
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 0.3922183333333333
$data[0,13] = 1.176655
$data[0,14] = 0.06257714082953221
$data[0,15] = 0.06257714082953222
$data[0,16] = 0.04435518688
$data[0,17] = 0.39919668192
$data[0,18] = 0.02045857591572871
$data[0,19] = 0.02045857591572871
$data[1,0] = "FAPs"
$data[1,1] = "Ntf5"
$data[1,2] = "Ngfr"
$data[1,3] = "MuSCs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.113088
$data[1,7] = 0.339264
$data[1,8] = 0.3269336956678857
$data[1,9] = 0.3269336956678857
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 5.787580333333334
$data[1,13] = 17.362741
$data[1,14] = 0.9233893441524432
$data[1,15] = 0.9233893441524432
$data[1,16] = 0.654505884736
$data[1,17] = 5.890552962624001
$data[1,18] = 0.3018870908241035
$data[1,19] = 0.3018870908241035
$data[2,0] = "FAPs"
$data[2,1] = "Ntf5"
$data[2,2] = "Ngfr"
$data[2,3] = "Neutrophils"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.113088
$data[2,7] = 0.339264
$data[2,8] = 0.3269336956678857
$data[2,9] = 0.3269336956678857
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.034619
$data[2,13] = 0.103857
$data[2,14] = 0.005523347213187152
$data[2,15] = 0.005523347213187152
$data[2,16] = 0.003914993472000001
$data[2,17] = 0.035234941248
$data[2,18] = 0.001805768316864193
$data[2,19] = 0.001805768316864193
$data[3,0] = "FAPs"
$data[3,1] = "Ntf5"
$data[3,2] = "Ngfr"
$data[3,3] = "Resolving-Mac"
$data[3,4] = 1
$data[3,5] = 0.3333333333333333
$data[3,6] = 0.113088
$data[3,7] = 0.339264
$data[3,8] = 0.3269336956678857
$data[3,9] = 0.3269336956678857
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.05333966666666667
$data[3,13] = 0.160019
$data[3,14] = 0.008510167804837369
$data[3,15] = 0.008510167804837371
$data[3,16] = 0.006032076224000001
$data[3,17] = 0.054288686016
$data[3,18] = 0.00278226061118934
$data[3,19] = 0.00278226061118934
$data[4,0] = "Inflammatory-Mac"
$data[4,1] = "Ntf5"
$data[4,2] = "Ngfr"
$data[4,3] = "ECs"
$data[4,4] = 1
$data[4,5] = 0.3333333333333333
$data[4,6] = 0.084206
$data[4,7] = 0.252618
$data[4,8] = 0.2434367817753429
$data[4,9] = 0.243436781775343
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 0.3922183333333333
$data[4,13] = 1.176655
$data[4,14] = 0.06257714082953221
$data[4,15] = 0.06257714082953222
$data[4,16] = 0.03302713697666667
$data[4,17] = 0.29724423279
$data[4,18] = 0.01523357777624374
$data[4,19] = 0.01523357777624374
$data[5,0] = "Inflammatory-Mac"
$data[5,1] = "Ntf5"
$data[5,2] = "Ngfr"
$data[5,3] = "MuSCs"
$data[5,4] = 1
$data[5,5] = 0.3333333333333333
$data[5,6] = 0.084206
$data[5,7] = 0.252618
$data[5,8] = 0.2434367817753429
$data[5,9] = 0.243436781775343
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 5.787580333333334
$data[5,13] = 17.362741
$data[5,14] = 0.9233893441524432
$data[5,15] = 0.9233893441524432
$data[5,16] = 0.4873489895486667
$data[5,17] = 4.386140905938
$data[5,18] = 0.2247869302661154
$data[5,19] = 0.2247869302661154
$data[6,0] = "Inflammatory-Mac"
$data[6,1] = "Ntf5"
$data[6,2] = "Ngfr"
$data[6,3] = "Neutrophils"
$data[6,4] = 1
$data[6,5] = 0.3333333333333333
$data[6,6] = 0.084206
$data[6,7] = 0.252618
$data[6,8] = 0.2434367817753429
$data[6,9] = 0.243436781775343
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.034619
$data[6,13] = 0.103857
$data[6,14] = 0.005523347213187152
$data[6,15] = 0.005523347213187152
$data[6,16] = 0.002915127514000001
$data[6,17] = 0.026236147626
$data[6,18] = 0.001344585870206089
$data[6,19] = 0.001344585870206089
$data[7,0] = "Inflammatory-Mac"
$data[7,1] = "Ntf5"
$data[7,2] = "Ngfr"
$data[7,3] = "Resolving-Mac"
$data[7,4] = 1
$data[7,5] = 0.3333333333333333
$data[7,6] = 0.084206
$data[7,7] = 0.252618
$data[7,8] = 0.2434367817753429
$data[7,9] = 0.243436781775343
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.05333966666666667
$data[7,13] = 0.160019
$data[7,14] = 0.008510167804837369
$data[7,15] = 0.008510167804837371
$data[7,16] = 0.004491519971333333
$data[7,17] = 0.040423679742
$data[7,18] = 0.002071687862777744
$data[7,19] = 0.002071687862777745
$data[8,0] = "MuSCs"
$data[8,1] = "Ntf5"
$data[8,2] = "Ngfr"
$data[8,3] = "ECs"
$data[8,4] = 1
$data[8,5] = 0.3333333333333333
$data[8,6] = 0.08859233333333333
$data[8,7] = 0.265777
$data[8,8] = 0.2561175274521424
$data[8,9] = 0.2561175274521424
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 0.3922183333333333
$data[8,13] = 1.176655
$data[8,14] = 0.06257714082953221
$data[8,15] = 0.06257714082953222
$data[8,16] = 0.03474753732611111
$data[8,17] = 0.312727835935
$data[8,18] = 0.0160271025842843
$data[8,19] = 0.0160271025842843
$data[9,0] = "MuSCs"
$data[9,1] = "Ntf5"
$data[9,2] = "Ngfr"
$data[9,3] = "MuSCs"
$data[9,4] = 1
$data[9,5] = 0.3333333333333333
$data[9,6] = 0.08859233333333333
$data[9,7] = 0.265777
$data[9,8] = 0.2561175274521424
$data[9,9] = 0.2561175274521424
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 5.787580333333334
$data[9,13] = 17.362741
$data[9,14] = 0.9233893441524432
$data[9,15] = 0.9233893441524432
$data[9,16] = 0.5127352460841111
$data[9,17] = 4.614617214757
$data[9,18] = 0.2364961956999792
$data[9,19] = 0.2364961956999792
$data[10,0] = "MuSCs"
$data[10,1] = "Ntf5"
$data[10,2] = "Ngfr"
$data[10,3] = "Neutrophils"
$data[10,4] = 1
$data[10,5] = 0.3333333333333333
$data[10,6] = 0.08859233333333333
$data[10,7] = 0.265777
$data[10,8] = 0.2561175274521424
$data[10,9] = 0.2561175274521424
$data[10,10] = 1
$data[10,11] = 0.3333333333333333
$data[10,12] = 0.034619
$data[10,13] = 0.103857
$data[10,14] = 0.005523347213187152
$data[10,15] = 0.005523347213187152
$data[10,16] = 0.003066977987666667
$data[10,17] = 0.027602801889
$data[10,18] = 0.001414626031501175
$data[10,19] = 0.001414626031501175
$data[11,0] = "MuSCs"
$data[11,1] = "Ntf5"
$data[11,2] = "Ngfr"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 1
$data[11,5] = 0.3333333333333333
$data[11,6] = 0.08859233333333333
$data[11,7] = 0.265777
$data[11,8] = 0.2561175274521424
$data[11,9] = 0.2561175274521424
$data[11,10] = 1
$data[11,11] = 0.3333333333333333
$data[11,12] = 0.05333966666666667
$data[11,13] = 0.160019
$data[11,14] = 0.008510167804837369
$data[11,15] = 0.008510167804837371
$data[11,16] = 0.004725485529222222
$data[11,17] = 0.04252936976299999
$data[11,18] = 0.002179603136377774
$data[11,19] = 0.002179603136377774
$data[12,0] = "Resolving-Mac"
$data[12,1] = "Ntf5"
$data[12,2] = "Ngfr"
$data[12,3] = "ECs"
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.06001866666666666
$data[12,7] = 0.180056
$data[12,8] = 0.1735119951046289
$data[12,9] = 0.1735119951046289
$data[12,10] = 2
$data[12,11] = 0.6666666666666666
$data[12,12] = 0.3922183333333333
$data[12,13] = 1.176655
$data[12,14] = 0.06257714082953221
$data[12,15] = 0.06257714082953222
$data[12,16] = 0.02354042140888889
$data[12,17] = 0.21186379268
$data[12,18] = 0.01085788455327547
$data[12,19] = 0.01085788455327547
$data[13,0] = "Resolving-Mac"
$data[13,1] = "Ntf5"
$data[13,2] = "Ngfr"
$data[13,3] = "MuSCs"
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.06001866666666666
$data[13,7] = 0.180056
$data[13,8] = 0.1735119951046289
$data[13,9] = 0.1735119951046289
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 5.787580333333334
$data[13,13] = 17.362741
$data[13,14] = 0.9233893441524432
$data[13,15] = 0.9233893441524432
$data[13,16] = 0.3473628548328889
$data[13,17] = 3.126265693496
$data[13,18] = 0.1602191273622453
$data[13,19] = 0.1602191273622453
$data[14,0] = "Resolving-Mac"
$data[14,1] = "Ntf5"
$data[14,2] = "Ngfr"
$data[14,3] = "Neutrophils"
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.06001866666666666
$data[14,7] = 0.180056
$data[14,8] = 0.1735119951046289
$data[14,9] = 0.1735119951046289
$data[14,10] = 1
$data[14,11] = 0.3333333333333333
$data[14,12] = 0.034619
$data[14,13] = 0.103857
$data[14,14] = 0.005523347213187152
$data[14,15] = 0.005523347213187152
$data[14,16] = 0.002077786221333333
$data[14,17] = 0.018700075992
$data[14,18] = 0.0009583669946156948
$data[14,19] = 0.0009583669946156948
$data[15,0] = "Resolving-Mac"
$data[15,1] = "Ntf5"
$data[15,2] = "Ngfr"
$data[15,3] = "Resolving-Mac"
$data[15,4] = 1
$data[15,5] = 0.3333333333333333
$data[15,6] = 0.06001866666666666
$data[15,7] = 0.180056
$data[15,8] = 0.1735119951046289
$data[15,9] = 0.1735119951046289
$data[15,10] = 1
$data[15,11] = 0.3333333333333333
$data[15,12] = 0.05333966666666667
$data[15,13] = 0.160019
$data[15,14] = 0.008510167804837369
$data[15,15] = 0.008510167804837371
$data[15,16] = 0.003201375673777778
$data[15,17] = 0.028812381064
$data[15,18] = 0.001476616194492512
$data[15,19] = 0.001476616194492513

$rng = $ws.Range("A2:T17")
$rng.Value = $data
